$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 119.6201171875
    3  = 126.5068817138672
    4  = 123.7748718261719
    5  = 119.1963729858398
    6  = 120.5098114013672
    7  = 119.1285552978516
    8  = 120.1011199951172
    9  = 117.042350769043
    10 = 119.5248489379883
    11 = 118.9073867797852
    12 = 120.6032791137695
    13 = 128.7219543457031
    14 = 132.1677551269531
    15 = 140.4591522216797
    16 = 163.1872863769531
    17 = 194.1389007568359
    18 = 184.5361938476562
    19 = 196.8851928710938
    20 = 200.9994506835938
    21 = 202.9446411132812
    22 = 203.6833953857422
    23 = 197.2630157470703
    24 = 199.1311798095703
    25 = 198.321044921875
    26 = 194.4628601074219
    27 = 196.814697265625
    28 = 193.0293884277344
    29 = 186.9048309326172
    30 = 189.3831024169922
    31 = 193.3955230712891
    32 = 199.1823577880859
    33 = 224.8252258300781
    34 = 216.3576507568359
    35 = 253.4040374755859
    36 = 251.6214904785156
    37 = 255.8715667724609
    38 = 228.4771270751953
    39 = 214.3982543945312
    40 = 204.7801818847656
    41 = 187.9519958496094
    42 = 177.8160858154297
    43 = 157.0883026123047
    44 = 160.8988342285156
    45 = 147.4034729003906
    46 = 150.0952453613281
    47 = 137.8648223876953
    48 = 145.6342163085938
    49 = 144.98779296875
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
